# Update the "sales/views" count column (F) on sheet "展览" and sheet "全部类型".
# Both sheets list the same events; the row numbering differs slightly because
# "全部类型" has one extra row inserted early in the sheet, so each sheet's row
# map is specified separately.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet4 = $wb.Worksheets.Item("全部类型")

# Row => new value, for sheet "展览"
$updatesSheet1 = @{
    3  = 1891
    4  = 157
    6  = 3197
    7  = 587
    9  = 294
    10 = 660
    11 = 556
    12 = 562
    13 = 411
    14 = 145
    15 = 1801
    16 = 1394
    18 = 1656
    21 = 624
    23 = 49
    30 = 47
    32 = 4051
    33 = 16
    34 = 785
    36 = 1747
    38 = 1919
}

# Row => new value, for sheet "全部类型"
$updatesSheet4 = @{
    3  = 1891
    4  = 157
    6  = 3197
    7  = 587
    9  = 294
    10 = 660
    11 = 556
    12 = 562
    14 = 411
    15 = 145
    16 = 1801
    17 = 1394
    19 = 1656
    22 = 624
    24 = 49
    31 = 47
    33 = 4051
    35 = 16
    37 = 785
    39 = 1747
    41 = 1919
}

foreach ($row in $updatesSheet1.Keys) {
    $sheet1.Range("F$row").Value = $updatesSheet1[$row]
}

foreach ($row in $updatesSheet4.Keys) {
    $sheet4.Range("F$row").Value = $updatesSheet4[$row]
}
